# Generate Report for Handoff
#
# A fresh handoff-report run moved the zh-cn / de-de localization jobs from
# "In Translation" to "Ready for handoff" and refreshed the associated
# timestamps on all three sheets (Overview, zh-cn, de-de). Excel then
# auto-sized the Status columns to fit the new, longer status text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: zh-cn / de-de status + HO xliff generate date ------
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-10-21 04:31:22"

# --- zh-cn sheet: Status + Latest Handoff Datetime -----------------------
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-10-21 04:31:11"

# --- de-de sheet: Status + Latest Handoff Datetime -----------------------
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-10-21 04:31:22"

# --- Widen the Status columns to fit "Ready for handoff" -----------------
# (Overview columns E & F are the zh-cn / de-de status columns; the zh-cn
# and de-de sheets' own Status column is column C.)
$overview.Columns.Item(5).ColumnWidth = 16.3826548258464
$overview.Columns.Item(6).ColumnWidth = 16.3826548258464
$zhcn.Columns.Item(3).ColumnWidth = 16.3826548258464
$dede.Columns.Item(3).ColumnWidth = 16.3826548258464
